$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.814.08"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.804.90"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'603.95"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "'166.10"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11").Value = "'6.32"
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "'0.0000252"
$ws.Range("E12").Value = "  -0.82%  "
$ws.Range("D13").Value = "'36.03"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "4.444.82"
$ws.Range("E14").Value = "  +0.74%  "
$ws.Range("D15").Value = "3.818.84"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "67.834.28"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D20").Value = "'464.18"
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'12.14"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").Value = "'10.03"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "3.954.17"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'29.47"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").Value = "'9.10"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "'5.82"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D43").Value = "'44.18"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "'27.99"
$ws.Range("E46").Value = "  +9.01%  "
$ws.Range("D47").Value = "'150.83"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("D48").Value = "'1.39"
$ws.Range("E48").Value = "  +11.84%  "
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").Value = "'390.70"
$ws.Range("E51").Value = "  +0.21%  "
